# Adiciona um novo diapositivo (3) com o esquema "Título e objecto"
# (título + marcador de posição de conteúdo), tal como o PowerPoint faz
# quando se usa Base Diapositivo > Novo Diapositivo com esse esquema.

$p = $ppt.ActivePresentation

# Esquema 2 = ppLayoutText ("Título e Objecto") -> título + conteúdo,
# que é o layout já usado no diapositivo 2 (slideLayout2.xml).
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)
